$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "58.607.57"
Set-TextValue "E2" "  +1.79%  "

Set-TextValue "D3" "3.164.32"
Set-TextValue "E3" "  +1.73%  "

Set-TextValue "E4" "  +0.01%  "

Set-TextValue "D5" "529.57"
Set-TextValue "E5" "  -0.21%  "

Set-TextValue "D6" "139.80"
Set-TextValue "E6" "  +1.19%  "

Set-TextValue "E7" "  +0.07%  "

Set-TextValue "D8" "0.540"
Set-TextValue "E8" "  +14.78%  "

Set-TextValue "D9" "7.30"
Set-TextValue "E9" "  -0.08%  "

Set-TextValue "D10" "0.437"
Set-TextValue "E10" "  +5.86%  "

Set-TextValue "E11" "  +3.95%  "

Set-TextValue "E12" "  +2.56%  "

Set-TextValue "D13" "3.710.17"
Set-TextValue "E13" "  +1.97%  "

Set-TextValue "D14" "25.72"
Set-TextValue "E14" "  +0.14%  "

Set-TextValue "E15" "  +3.65%  "

Set-TextValue "D16" "58.663.55"
Set-TextValue "E16" "  +1.69%  "

Set-TextValue "D17" "6.25"
Set-TextValue "E17" "  +3.68%  "

Set-TextValue "D18" "3.108.43"
Set-TextValue "E18" "  +0.30%  "

Set-TextValue "D19" "12.99"
Set-TextValue "E19" "  +2.51%  "

Set-TextValue "B20" "Uniswap"
Set-TextValue "C20" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D20" "8.11"
Set-TextValue "E20" "  +0.37%  "

Set-TextValue "B21" "BitcoinCash"
Set-TextValue "C21" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D21" "375.83"
Set-TextValue "E21" "  +4.26%  "

Set-TextValue "E22" "  +0.22%  "

Set-TextValue "E23" "  +5.24%  "

Set-TextValue "D24" "69.70"
Set-TextValue "E24" "  +1.27%  "

Set-TextValue "E25" "  +0.50%  "

Set-TextValue "E26" "  +0.18%  "

Set-TextValue "D27" "8.26"
Set-TextValue "E27" "  +12.72%  "

Set-TextValue "D28" "0.0₃0866"
Set-TextValue "E28" "  -0.05%  "

Set-TextValue "D29" "22.36"
Set-TextValue "E29" "  +4.69%  "

Set-TextValue "D30" "1.88"
Set-TextValue "E30" "  +0.87%  "

Set-TextValue "E31" "  -0.65%  "

Set-TextValue "E32" "  +2.09%  "

Set-TextValue "E33" "  +1.18%  "

Set-TextValue "D34" "6.31"
Set-TextValue "E34" "  +4.31%  "

Set-TextValue "D35" "156.61"
Set-TextValue "E35" "  -1.71%  "

Set-TextValue "E36" "  +4.68%  "

Set-TextValue "D37" "25.06"
Set-TextValue "E37" "  -1.91%  "

Set-TextValue "D38" "2.677.24"
Set-TextValue "E38" "  +7.44%  "

Set-TextValue "D39" "0.0695"
Set-TextValue "E39" "  +3.90%  "

Set-TextValue "E40" "  +1.88%  "

Set-TextValue "E41" "  +6.86%  "

Set-TextValue "E42" "  +3.89%  "

Set-TextValue "E43" "  +4.06%  "

Set-TextValue "D44" "0.0288"
Set-TextValue "E44" "  +6.95%  "

Set-TextValue "E45" "  -0.03%  "

Set-TextValue "D46" "3.207.34"
Set-TextValue "E46" "  +1.85%  "

Set-TextValue "E47" "  +14.22%  "

Set-TextValue "D48" "6.23"
Set-TextValue "E48" "  +2.61%  "

Set-TextValue "D49" "0.981"
Set-TextValue "E49" "  -0.50%  "

Set-TextValue "D50" "20.10"
Set-TextValue "E50" "  +1.77%  "

Set-TextValue "E51" "  +1.81%  "
